# Update December 2024 statistics
# Fills in Circulation / ILL Loans / ILL Borrows figures (columns B, C, D)
# for each library row (rows 3-59) on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Each entry: Row, ColumnB (Circulation), ColumnC (ILL Loans), ColumnD (ILL Borrows)
$data = @(
    @(3, 70513, 10922, 11835),
    @(4, 41652, 3882, 5218),
    @(5, 110588, 10614, 9894),
    @(6, 2019, 1094, 303),
    @(7, 73874, 12756, 9124),
    @(8, 8519, 1699, 1817),
    @(9, 8829, 1878, 1260),
    @(10, 4154, 647, 340),
    @(11, 1551, 475, 4),
    @(12, 0, 0, 0),
    @(13, 1428, 304, 448),
    @(14, 4064, 1539, 1477),
    @(15, 7553, 2588, 1542),
    @(16, 5040, 2778, 618),
    @(17, 4121, 1351, 535),
    @(18, 27308, 3914, 5174),
    @(19, 2129, 992, 546),
    @(20, 27545, 3866, 4293),
    @(21, 557, 627, 185),
    @(22, 27424, 3099, 4782),
    @(23, 2005, 1099, 335),
    @(24, 30629, 3090, 5685),
    @(25, 122510, 10493, 13846),
    @(26, 9381, 3716, 1335),
    @(27, 0, 0, 0),
    @(28, 8233, 1898, 2033),
    @(29, 3782, 863, 789),
    @(30, 23362, 4144, 4695),
    @(31, 707, 168, 435),
    @(32, 4490, 2568, 477),
    @(33, 22008, 5266, 4470),
    @(34, 17455, 5030, 3214),
    @(35, 7884, 1011, 1623),
    @(36, 87796, 10138, 8585),
    @(37, 12404, 4635, 1668),
    @(38, 42814, 3038, 4102),
    @(39, 1881, 1660, 308),
    @(40, 2161, 811, 859),
    @(41, 4060, 874, 178),
    @(42, 17934, 866, 505),
    @(43, 414, 303, 77),
    @(44, 1329, 250, 96),
    @(45, 0, 0, 0),
    @(46, 5182, 1683, 645),
    @(47, 21868, 5312, 4178),
    @(48, 51673, 5320, 7745),
    @(49, 23914, 5318, 2165),
    @(50, 19597, 2384, 4318),
    @(51, 51224, 4867, 7995),
    @(52, 8161, 1564, 1974),
    @(53, 16889, 3361, 2833),
    @(54, 3537, 2075, 1440),
    @(55, 3342, 2146, 245),
    @(56, 8461, 1714, 3383),
    @(57, 18085, 7209, 4085),
    @(58, 23737, 2293, 856),
    @(59, 1065955, 162690, 151631)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
}
